# Reference-style edit: give all headings a uniform accent-1 color and
# rebalance the Heading2-4 weight/size ladder (see commit message).
$d = $word.ActiveDocument
$styles = $d.Styles

# Heading 1: drop the 65% "B5" theme-shade tint so it resolves to the
# same plain accent1 color the other heading levels already use.
$styles.Item("Heading1").Font.TextColor.ObjectThemeColor = 4

# Heading 2: 16pt -> 14pt (sz/szCs both carry the point size *2).
$styles.Item("Heading2").Font.Size = 14
$styles.Item("Heading2").Font.SizeBi = 14

# Heading 3: 14pt -> 12pt.
$styles.Item("Heading3").Font.Size = 12
$styles.Item("Heading3").Font.SizeBi = 12

# Heading 4: bold -> italic.
$styles.Item("Heading4").Font.Italic = $true
$styles.Item("Heading4").Font.Bold = $false

# Heading 5: no longer italic.
$styles.Item("Heading5").Font.Italic = $false
